$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-14 Wednesday" "2026-01-15 Thursday"

Replace-Text "237÷9=" "811÷4="
Replace-Text "884÷6=" "159÷5="
Replace-Text "588÷6=" "426÷2="
Replace-Text "733÷3=" "527÷8="
Replace-Text "684÷8=" "846÷6="
Replace-Text "591÷9=" "914÷9="
Replace-Text "475÷9=" "858÷2="
Replace-Text "922÷7=" "585÷9="
Replace-Text "842÷3=" "803÷5="
Replace-Text "826÷4=" "532÷2="
Replace-Text "572÷6=" "535÷4="
Replace-Text "526÷4=" "163÷5="
Replace-Text "667÷9=" "447÷3="
Replace-Text "749÷6=" "812÷6="
Replace-Text "249÷7=" "653÷9="
Replace-Text "503÷8=" "145÷2="
Replace-Text "128÷3=" "356÷9="
Replace-Text "123÷9=" "710÷3="
Replace-Text "684÷7=" "314÷5="
Replace-Text "443÷3=" "378÷9="
Replace-Text "131÷8=" "349÷9="
Replace-Text "949÷6=" "891÷6="
Replace-Text "261÷3=" "645÷6="
Replace-Text "929÷6=" "934÷9="
Replace-Text "777÷8=" "394÷4="
